# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns,
# and restore two coin rows that swapped ranking position (24/25, 37/38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.016.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "'2.415.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'554.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").Value = "'142.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "'2.412.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "'5.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").Value = "'26.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.93%  "
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.24%  "
$ws.Range("D16").Value = "'2.856.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Value = "'62.005.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "'2.417.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("D19").Value = "'11.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.81%  "
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "'323.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").Value = "'1.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'64.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").Value = "'9.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.76%  "
$ws.Range("D27").Value = "'576.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.33%  "
$ws.Range("D28").Value = "'2.536.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +4.13%  "
$ws.Range("D31").Value = "'0.0₃0928"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.26%  "
$ws.Range("D32").Value = "'1.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.76%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'5.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.77%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.89%  "
$ws.Range("D39").Value = "'0.384"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("E41").Value = "  +2.79%  "
$ws.Range("D42").Value = "'148.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'41.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("E45").Value = "  +6.09%  "
$ws.Range("D46").Value = "'2.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.85%  "
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("E48").Value = "  +5.40%  "
$ws.Range("D49").Value = "'20.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.18%  "
$ws.Range("E50").Value = "  +3.22%  "
$ws.Range("E51").Value = "  +1.61%  "
